{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst edits = [\n  { row: 0, col: 0, oldText: \"555\u00f79=61, 6\", newText: \"766\u00f78=95, 6\" },\n  { row: 0, col: 1, oldText: \"666\u00f73=222, 0\", newText: \"589\u00f74=147, 1\" },\n  { row: 0, col: 2, oldText: \"324\u00f77=46, 2\", newText: \"697\u00f73=232, 1\" },\n  { row: 0, col: 3, oldText: \"474\u00f73=158, 0\", newText: \"180\u00f72=90, 0\" },\n  { row: 0, col: 4, oldText: \"461\u00f75=92, 1\", newText: \"232\u00f76=38, 4\" },\n  { row: 4, col: 0, oldText: \"196\u00f79=21, 7\", newText: \"481\u00f77=68, 5\" },\n  { row: 4, col: 1, oldText: \"259\u00f75=51, 4\", newText: \"831\u00f74=207, 3\" },\n  { row: 4, col: 2, oldText: \"145\u00f74=36, 1\", newText: \"133\u00f73=44, 1\" },\n  { row: 4, col: 3, oldText: \"680\u00f73=226, 2\", newText: \"175\u00f77=25, 0\" },\n  { row: 4, col: 4, oldText: \"976\u00f73=325, 1\", newText: \"214\u00f72=107, 0\" },\n  { row: 8, col: 0, oldText: \"889\u00f75=177, 4\", newText: \"133\u00f79=14, 7\" },\n  { row: 8, col: 1, oldText: \"870\u00f79=96, 6\", newText: \"229\u00f78=28, 5\" },\n  { row: 8, col: 2, oldText: \"403\u00f74=100, 3\", newText: \"578\u00f77=82, 4\" },\n  { row: 8, col: 3, oldText: \"158\u00f76=26, 2\", newText: \"255\u00f76=42, 3\" },\n  { row: 8, col: 4, oldText: \"155\u00f78=19, 3\", newText: \"357\u00f73=119, 0\" },\n  { row: 12, col: 0, oldText: \"198\u00f79=22, 0\", newText: \"181\u00f77=25, 6\" },\n  { row: 12, col: 1, oldText: \"877\u00f73=292, 1\", newText: \"957\u00f73=319, 0\" },\n  { row: 12, col: 2, oldText: \"309\u00f78=38, 5\", newText: \"479\u00f77=68, 3\" },\n  { row: 12, col: 3, oldText: \"166\u00f73=55, 1\", newText: \"113\u00f75=22, 3\" },\n  { row: 12, col: 4, oldText: \"964\u00f77=137, 5\", newText: \"522\u00f78=65, 2\" },\n  { row: 16, col: 0, oldText: \"133\u00f73=44, 1\", newText: \"436\u00f75=87, 1\" },\n  { row: 16, col: 1, oldText: \"540\u00f77=77, 1\", newText: \"720\u00f79=80, 0\" },\n  { row: 16, col: 2, oldText: \"945\u00f78=118, 1\", newText: \"197\u00f79=21, 8\" },\n  { row: 16, col: 3, oldText: \"497\u00f79=55, 2\", newText: \"695\u00f73=231, 2\" },\n  { row: 16, col: 4, oldText: \"124\u00f78=15, 4\", newText: \"303\u00f75=60, 3\" },\n];\n\n// Load each target cell's first paragraph so we can replace its range in place\n// (this preserves paragraph/run formatting, unlike body.insertText(..., replace)).\nconst cellParas = [];\nfor (const edit of edits) {\n  const cell = table.getCellOrNullObject(edit.row, edit.col);\n  cell.load(\"body\");\n  cellParas.push({ edit, cell });\n}\nawait context.sync();\n\nconst paras = [];\nfor (const { edit, cell } of cellParas) {\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items\");\n  paras.push({ edit, cell, paragraphs });\n}\nawait context.sync();\n\nfor (const { edit, paragraphs } of paras) {\n  const para = paragraphs.items[0];\n  const range = para.getRange();\n  range.insertText(edit.newText, Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# (row, col, newText) triples \u2014 1-indexed Table.Cell() coordinates matching the\n# five populated data rows (table rows 1, 5, 9, 13, 17) x five columns each.\n$edits = @(\n    @{ Row = 1; Col = 1; NewText = '766\u00f78=95, 6' }\n    @{ Row = 1; Col = 2; NewText = '589\u00f74=147, 1' }\n    @{ Row = 1; Col = 3; NewText = '697\u00f73=232, 1' }\n    @{ Row = 1; Col = 4; NewText = '180\u00f72=90, 0' }\n    @{ Row = 1; Col = 5; NewText = '232\u00f76=38, 4' }\n    @{ Row = 5; Col = 1; NewText = '481\u00f77=68, 5' }\n    @{ Row = 5; Col = 2; NewText = '831\u00f74=207, 3' }\n    @{ Row = 5; Col = 3; NewText = '133\u00f73=44, 1' }\n    @{ Row = 5; Col = 4; NewText = '175\u00f77=25, 0' }\n    @{ Row = 5; Col = 5; NewText = '214\u00f72=107, 0' }\n    @{ Row = 9; Col = 1; NewText = '133\u00f79=14, 7' }\n    @{ Row = 9; Col = 2; NewText = '229\u00f78=28, 5' }\n    @{ Row = 9; Col = 3; NewText = '578\u00f77=82, 4' }\n    @{ Row = 9; Col = 4; NewText = '255\u00f76=42, 3' }\n    @{ Row = 9; Col = 5; NewText = '357\u00f73=119, 0' }\n    @{ Row = 13; Col = 1; NewText = '181\u00f77=25, 6' }\n    @{ Row = 13; Col = 2; NewText = '957\u00f73=319, 0' }\n    @{ Row = 13; Col = 3; NewText = '479\u00f77=68, 3' }\n    @{ Row = 13; Col = 4; NewText = '113\u00f75=22, 3' }\n    @{ Row = 13; Col = 5; NewText = '522\u00f78=65, 2' }\n    @{ Row = 17; Col = 1; NewText = '436\u00f75=87, 1' }\n    @{ Row = 17; Col = 2; NewText = '720\u00f79=80, 0' }\n    @{ Row = 17; Col = 3; NewText = '197\u00f79=21, 8' }\n    @{ Row = 17; Col = 4; NewText = '695\u00f73=231, 2' }\n    @{ Row = 17; Col = 5; NewText = '303\u00f75=60, 3' }\n)\n\nforeach ($edit in $edits) {\n    $cell = $t.Cell($edit.Row, $edit.Col)\n    $rng = $cell.Range\n    # Trim the trailing cell-mark (and any paragraph mark) so only the visible\n    # text is replaced; this preserves the run/paragraph formatting in place.\n    $rng.End = $rng.End - 1\n    $rng.Text = $edit.NewText\n}"}
